# Generate Report for Handback
# - Updates the "Status" text (for the 68f3f086... file) from
#   "Ready for handoff" to "Handback transform failed" on the Overview,
#   zh-cn and de-de sheets.
# - Records the handback/handoff filename mismatch error detail for that
#   file on both the zh-cn and de-de per-language sheets.
# - Widens the "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: row 3 is the 68f3f086-e976-4ec2-afbc-238b289183b0.md file,
# columns E (zh-cn) and F (de-de) hold its per-language status.
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# zh-cn / de-de sheets: row 3 is the same file; column C is "Status".
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# Error Detail (column P) messages describing the handback/handoff file
# name mismatch detected for this file in each language.
$zhcn.Range("P3").Value = "Handback file name: mn0nenbu.zq0 is different with handoff file name: 68f3f086-e976-4ec2-afbc-238b289183b0.6a49ebd4c3fbc5fa854a18954d3f92595a411e3d.zh-cn."
$dede.Range("P3").Value = "Handback file name: mn0nenbu.zq0 is different with handoff file name: 68f3f086-e976-4ec2-afbc-238b289183b0.6a49ebd4c3fbc5fa854a18954d3f92595a411e3d.de-de."

# Widen the Error Detail column (P = 16th column) on both language sheets
# so the new, longer messages are visible. (39.17 is the ColumnWidth that
# Excel persists to OOXML as width="40", matching column A's own width.)
$zhcn.Columns.Item(16).ColumnWidth = 39.17
$dede.Columns.Item(16).ColumnWidth = 39.17
